$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
  @("Sending cluster", "Ligand symbol", "Receptor symbol", "Target cluster", "Ligand-expressing cells", "Ligand detection rate", "Ligand average expression value", "Ligand total expression value", "Ligand derived specificity of average expression value", "Ligand derived specificity of total expression value", "Receptor-expressing cells", "Receptor detection rate", "Receptor average expression value", "Receptor total expression value", "Receptor derived specificity of average expression value", "Receptor derived specificity of total expression value", "Edge average expression weight", "Edge total expression weight", "Edge average expression derived specificity", "Edge total expression derived specificity"),
  @("ECs", "Ifnk", "Ifnar1", "ECs", 1, 0.3333333333333333, 0.2027413333333333, 0.608224, 0.1911300550679205, 0.1911300550679204, 3, 1, 27.30988266666667, 81.929648, 0.2168690090390243, 0.2168690090390242, 5.536842025016889, 49.831578225152, 0.04145018564015405, 0.04145018564015403),
  @("ECs", "Ifnk", "Ifnar1", "FAPs", 1, 0.3333333333333333, 0.2027413333333333, 0.608224, 0.1911300550679205, 0.1911300550679204, 3, 1, 20.371636, 61.114908, 0.1617720795708915, 0.1617720795708915, 4.130172644821333, 37.171553803392, 0.03091950647683651, 0.03091950647683651),
  @("ECs", "Ifnk", "Ifnar1", "Inflammatory-Mac", 1, 0.3333333333333333, 0.2027413333333333, 0.608224, 0.1911300550679205, 0.1911300550679204, 3, 1, 42.52135866666666, 127.564076, 0.3376640254953711, 0.3376640254953711, 8.620836951224888, 77.587532561024, 0.06453774378738597, 0.06453774378738596),
  @("ECs", "Ifnk", "Ifnar1", "MuSCs", 1, 0.3333333333333333, 0.2027413333333333, 0.608224, 0.1911300550679205, 0.1911300550679204, 3, 1, 6.417394333333333, 19.252183, 0.05096081761571768, 0.05096081761571768, 1.301071083665778, 11.709639752992, 0.009740143877198373, 0.00974014387719837),
  @("ECs", "Ifnk", "Ifnar1", "Resolving-Mac", 1, 0.3333333333333333, 0.2027413333333333, 0.608224, 0.1911300550679205, 0.1911300550679204, 3, 1, 29.307738, 87.923214, 0.2327340682789955, 0.2327340682789955, 5.941889879104, 53.477008911936, 0.04448247528634557, 0.04448247528634556),
  @("Inflammatory-Mac", "Ifnk", "Ifnar1", "ECs", 1, 0.3333333333333333, 0.7050673333333334, 2.115202, 0.6646871460839683, 0.6646871460839682, 3, 1, 27.30988266666667, 81.929648, 0.2168690090390243, 0.2168690090390242, 19.25530614543289, 173.297755308896, 0.1441500426922074, 0.1441500426922073),
  @("Inflammatory-Mac", "Ifnk", "Ifnar1", "FAPs", 1, 0.3333333333333333, 0.7050673333333334, 2.115202, 0.6646871460839683, 0.6646871460839682, 3, 1, 20.371636, 61.114908, 0.1617720795708915, 0.1617720795708915, 14.36337507015733, 129.270375631416, 0.1075278218860445, 0.1075278218860445),
  @("Inflammatory-Mac", "Ifnk", "Ifnar1", "Inflammatory-Mac", 1, 0.3333333333333333, 0.7050673333333334, 2.115202, 0.6646871460839683, 0.6646871460839682, 3, 1, 42.52135866666666, 127.564076, 0.3376640254953711, 0.3376640254953711, 29.98042096481689, 269.823788683352, 0.2244409374417425, 0.2244409374417425),
  @("Inflammatory-Mac", "Ifnk", "Ifnar1", "MuSCs", 1, 0.3333333333333333, 0.7050673333333334, 2.115202, 0.6646871460839683, 0.6646871460839682, 3, 1, 6.417394333333333, 19.252183, 0.05096081761571768, 0.05096081761571768, 4.524695109551778, 40.722255985966, 0.033873000423097, 0.033873000423097),
  @("Inflammatory-Mac", "Ifnk", "Ifnar1", "Resolving-Mac", 1, 0.3333333333333333, 0.7050673333333334, 2.115202, 0.6646871460839683, 0.6646871460839682, 3, 1, 29.307738, 87.923214, 0.2327340682789955, 0.2327340682789955, 20.663928677692, 185.975358099228, 0.1546953436408769, 0.1546953436408769),
  @("Resolving-Mac", "Ifnk", "Ifnar1", "ECs", 1, 0.3333333333333333, 0.152942, 0.458826, 0.1441827988481113, 0.1441827988481113, 3, 1, 27.30988266666667, 81.929648, 0.2168690090390243, 0.2168690090390242, 4.176828074805333, 37.591452673248, 0.03126878070666287, 0.03126878070666287),
  @("Resolving-Mac", "Ifnk", "Ifnar1", "FAPs", 1, 0.3333333333333333, 0.152942, 0.458826, 0.1441827988481113, 0.1441827988481113, 3, 1, 20.371636, 61.114908, 0.1617720795708915, 0.1617720795708915, 3.115678753112, 28.041108778008, 0.02332475120801052, 0.02332475120801052),
  @("Resolving-Mac", "Ifnk", "Ifnar1", "Inflammatory-Mac", 1, 0.3333333333333333, 0.152942, 0.458826, 0.1441827988481113, 0.1441827988481113, 3, 1, 42.52135866666666, 127.564076, 0.3376640254953711, 0.3376640254953711, 6.503301637197333, 58.529714734776, 0.04868534426624262, 0.04868534426624262),
  @("Resolving-Mac", "Ifnk", "Ifnar1", "MuSCs", 1, 0.3333333333333333, 0.152942, 0.458826, 0.1441827988481113, 0.1441827988481113, 3, 1, 6.417394333333333, 19.252183, 0.05096081761571768, 0.05096081761571768, 0.9814891241286666, 8.833402117158, 0.007347673315422311, 0.007347673315422311),
  @("Resolving-Mac", "Ifnk", "Ifnar1", "Resolving-Mac", 1, 0.3333333333333333, 0.152942, 0.458826, 0.1441827988481113, 0.1441827988481113, 3, 1, 29.307738, 87.923214, 0.2327340682789955, 0.2327340682789955, 4.482384065196, 40.341456586764, 0.03355624935177301, 0.03355624935177301)
)

for ($i = 0; $i -lt $data.Length; $i++) {
  $r = $data[$i]
  for ($j = 0; $j -lt $r.Length; $j++) {
    $ws.Cells.Item($i + 1, $j + 1).Value = $r[$j]
  }
}
